# Updated symbol list on Tue Dec 13 23:56:16 UTC 2022 with GitHub Actions
# Applies the price/volume/ranking refresh described in the commit diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds text-formatted numbers (e.g. "271.82"). A leading apostrophe
# forces Excel to store them as text rather than re-parsing them as numeric
# values (which would change cell type / introduce floating point noise).

$ws.Range("D2").Value = "'271.82"
$ws.Range("D3").Value = "'22.87"
$ws.Range("D4").Value = "'6.458"
$ws.Range("D5").Value = "'0.06225"
$ws.Range("D7").Value = "'6.664"
$ws.Range("D8").Value = "'1.385"
$ws.Range("D9").Value = "'0.8313"
$ws.Range("D10").Value = "'0.01381"
$ws.Range("D11").Value = "'0.1603"
$ws.Range("D12").Value = "'0.08250"
$ws.Range("D13").Value = "'0.03429"
$ws.Range("D14").Value = "'0.03186"
$ws.Range("D15").Value = "'0.09353"
$ws.Range("D16").Value = "'3.839"
$ws.Range("D17").Value = "'0.001647"
$ws.Range("D18").Value = "'0.04746"
$ws.Range("D19").Value = "'0.006330"
$ws.Range("D20").Value = "'0.005686"
$ws.Range("D21").Value = "'0.001076"
$ws.Range("D22").Value = "'0.0001500"
$ws.Range("D23").Value = "'3.713"
$ws.Range("D24").Value = "'2.400"
$ws.Range("D25").Value = "'0.3347"
$ws.Range("D26").Value = "'0.1252"
$ws.Range("D40").Value = "'0.04703"
$ws.Range("D41").Value = "'0.007052"
$ws.Range("D42").Value = "'0.003799"
$ws.Range("E42").Value = "41CEJICEJI"
$ws.Range("D43").Value = "'0.1160"
$ws.Range("D45").Value = "'0.00006282"
$ws.Range("D47").Value = "'0.00000000750"

# Rows 49/50: BOLO moves ahead of CryptobidCoin in the ranking, each with
# refreshed prices and "Volume(1h)" labels.
$ws.Range("B49").Value = "BOLO"
$ws.Range("C49").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$ws.Range("D49").Value = "'0.002140"
$ws.Range("E49").Value = "48BOLOBOLOBestin24h"

$ws.Range("B50").Value = "CryptobidCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/h39bvStAP+cryptobidcoin-cbc"
$ws.Range("D50").Value = "'0.00001400"
$ws.Range("E50").Value = "49CryptobidCoinCBCWorstin24h"

$ws.Range("D51").Value = "'0.01240"
